$d = $word.ActiveDocument

function Find-ParaIndexByText($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $i
        }
    }
    return -1
}

# Inserts a new paragraph right after paragraph index $afterIdx, applies
# $styleName (if given) and $text (if given) to it, and returns the index
# of the newly created paragraph.
function Insert-ParaAfter($doc, $afterIdx, $styleName, $text) {
    $anchor = $doc.Paragraphs.Item($afterIdx)
    $anchor.Range.InsertParagraphAfter() | Out-Null
    $newIdx = $afterIdx + 1
    $newPara = $doc.Paragraphs.Item($newIdx)
    if ($styleName) {
        $newPara.Style = $styleName
    }
    if ($text) {
        $newPara.Range.Text = $text
    }
    return $newIdx
}

# ---------------------------------------------------------------------
# 1) Footer date update (Last Updated ...)
# ---------------------------------------------------------------------
$sec = $d.Sections.First
$ftr = $sec.Footers.Item(1)
$ftr.Range.Find.Execute("12/04/2023 15:05", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "13/03/2024 02:05", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Version 0.0.4 -> Added: "Lock Mechanics." and "Note Reading."
#    (inserted right after "Hide mechanics.", reusing its bullet list)
# ---------------------------------------------------------------------
$idxHideMech = Find-ParaIndexByText $d "Hide mechanics."
$idxLock = Insert-ParaAfter $d $idxHideMech $null "Lock Mechanics."
$idxNote = Insert-ParaAfter $d $idxLock $null "Note Reading."

# ---------------------------------------------------------------------
# 3) Version 0.0.4 -> Fixed: "An issue where multiple notes could be opened."
#    (inserted right after the existing last "Fixed" bullet)
# ---------------------------------------------------------------------
$idxHideWhileHiding = Find-ParaIndexByText $d "An issue with the ability to open inventory while hiding."
Insert-ParaAfter $d $idxHideWhileHiding $null "An issue where multiple notes could be opened." | Out-Null

# ---------------------------------------------------------------------
# 4) New "Version 0.0.5" section at the end of the document
# ---------------------------------------------------------------------
$idxUI = Find-ParaIndexByText $d "UI and HUD to display Health/Stamina/Battery Life."

# blank paragraph right after "UI and HUD..."
$idxBlank1 = Insert-ParaAfter $d $idxUI $null $null

# "Version 0.0.5" (Heading 1)
$idxVersion = Insert-ParaAfter $d $idxBlank1 "Heading 1" "Version 0.0.5"

# "Added" (Heading 2)
$idxAdded = Insert-ParaAfter $d $idxVersion "Heading 2" "Added"

# "Jump Scare System" bullet (reuse the existing bullet-list definition
# already used by numId 7, e.g. the "Hide mechanics." bullet)
$existingBulletPara = $d.Paragraphs.Item((Find-ParaIndexByText $d "Hide mechanics."))
$bulletTemplate = $existingBulletPara.Range.ListFormat.ListTemplate

$idxJump = Insert-ParaAfter $d $idxAdded "List Paragraph" "Jump Scare System"
$pJump = $d.Paragraphs.Item($idxJump)
$pJump.Range.ListFormat.ApplyListTemplate($bulletTemplate)

# "Items" (Heading 3)
$idxItems = Insert-ParaAfter $d $idxJump "Heading 3" "Items"

# "Fixed" (Heading 2)
$idxFixed = Insert-ParaAfter $d $idxItems "Heading 2" "Fixed"

# "Updated" (Heading 2)
$idxUpdated = Insert-ParaAfter $d $idxFixed "Heading 2" "Updated"

# two trailing blank paragraphs
$idxBlank2 = Insert-ParaAfter $d $idxUpdated $null $null
Insert-ParaAfter $d $idxBlank2 $null $null | Out-Null
